$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("ZoneLetter") rows 2 through 467 currently hold "T";
# change each of them to "V".
$ws.Range("E2:E467").Value = "V"
